# Performance Metrics - Garrett PC.xlsx
# "Reworked cuda files to run kernel with 64 cores"
#
# Fills in the measured results for the Knapsack, KnapsackCuda and
# MatrixMultiply rows (rows 6, 7 and 8) that were previously blank, adds
# a formula to B4 (VectorAdd AvgExecTime) that was a bare literal value,
# nudges the default column width slightly, and moves the active
# selection to D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Selection moves from H5 to D6 -----------------------------------
[void]$ws.Range("D6").Select()

# --- Sheet-wide default column width (cosmetic, best effort) --------
$ws.StandardWidth = 11.625

# --- VectorAdd (row 4): literal value becomes its own formula -------
$ws.Range("B4").Formula = "=0.002838102"

# --- Knapsack (row 6) -------------------------------------------------
$ws.Range("B6").Formula = "=SUM(0.00353895,0.003344634,0.003759654,0.003101237,0.001890542,0.004512839,0.002438622,0.004772196,0.003817112,0.003118089)/10"
$ws.Range("C6").Formula = "=SUM(355423,376613,354202,399117,367524,381675,382102,407409,422486,360049)/10"

# --- KnapsackCuda (row 7) --------------------------------------------
$ws.Range("B7").Formula = "=0.0000201"
$ws.Range("B7").NumberFormat = "#,##0.000000000"
$ws.Range("B7").Font.Bold = $true
$ws.Range("F7").Value = 27449
$ws.Range("G7").Value = 3352
$ws.Range("H7").Formula = "=F7/G7"

# --- MatrixMultiply (row 8) ------------------------------------------
$ws.Range("B8").Formula = "=SUM(0.004614999,0.002712369,0.003579115,0.003510445,0.003024007,0.003475369,0.003888869,0.002700506,0.003084371)/10"
$ws.Range("C8").Formula = "=SUM(422681,343206,394682,390143,387199,398692,360587,38228,354,496,333309)/10"
$ws.Range("D8").Formula = "=SUM(119864,90868,94379,105399,103616,97594,97703,106168,92069,85130)/10"
$ws.Range("E8").Formula = "=SUM(136,130,134,134,132,134,134,133,133,132)/10"
$ws.Range("F8").Formula = "=SUM(46052033, 32908310, 32564908, 32358966, 23541511, 24745511, 30679930, 34695561, 30556292, 33679487)/10"
$ws.Range("G8").Formula = "=SUM(54227967, 41235861, 43352442, 43855952, 27684759, 28750882, 37685914, 39955331, 39385853, 39385853)/10"
$ws.Range("H8").Formula = "=F8/G8"

[void]$wb.Save()
